$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Тест-кейсы на автоматизацию")

# Row 14: change D14 from "To do" to "Automated"
$ws.Range("D14").Value = "Automated"

# Row 15: clear B15, C15, D15 (keep E15 as-is, it's already empty)
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()
$ws.Range("D15").ClearContents()

# Update active cell selection to D15
$ws.Range("D15").Select()
